# SW1116_noCTRL_meas.xlsx bug fix:
#  - Sheet1 had 43 stray rows (45:87) that only carried a leftover index
#    value in column A with no real measurement data next to them; delete
#    them so the sheet's used range shrinks back to A1:N44.
#  - Make Sheet1 the active/selected sheet & tab (it was Sheet3 before),
#    and leave the view scrolled/selected where the user last left it
#    (selection E58) after trimming the bogus rows.
#  - Sheet3 is no longer the active tab, so it loses its tabSelected flag;
#    its own selection (A2:N44) is left untouched.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Remove the extraneous rows 45-87 on Sheet1 (only column A had values,
# numbering 44..86 with no corresponding measurements) -> dimension
# becomes A1:N44.
$ws1.Rows("45:87").Delete()

# Sheet1 becomes the active sheet/tab instead of Sheet3.
$ws1.Activate()

# Restore the view state left on Sheet1: scrolled down with E58 selected.
$ws1.Range("E58").Select()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
